# Correct mistake in water characterization factors for exiobase water flows.
# The values in columns B, C, D, E (rows 4-101) were stored in the wrong
# unit/scale; multiply each numeric value by 1,000,000 to fix it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 4; $r -le 101; $r++) {
    foreach ($c in 2..5) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -ne $null) {
            $cell.Value = $v * 1000000
        }
    }
}

# Update the view: scroll position back to top-left and move the
# selection from A93 to F9.
[void]$ws.Range("F9").Select()

